# Append three new ticket rows (140-142) loaded from the admin JSON feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-05-21", "12:17:26", "Etiquetadora2", "-", "-", "-", "-", "12:17:28", "0:00:02"),
    @("2024-05-21", "12:22:15", "Etiquetadora2", "-", "-", "-", "-", "12:22:17", "0:00:02"),
    @("2024-05-21", "12:23:29", "Etiquetadora",  "-", "-", "-", "-", "12:23:31", "0:00:02")
)

$startRow = 140
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Column A holds a date string ("2024-05-21"); format it as Text first so
    # Excel stores the literal string instead of auto-converting it to a
    # date serial number.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]

    for ($c = 1; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
